$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Update existing column-B estimates for rows 278-304 (date unchanged)
# ---------------------------------------------------------------------
$updatedValues = @(
    @(278, 2.151317834854126),
    @(279, 2.263473272323608),
    @(280, 2.330948829650879),
    @(281, 2.471001148223877),
    @(282, 2.781925678253174),
    @(283, 2.807555198669434),
    @(284, 2.854128122329712),
    @(285, 2.642922163009644),
    @(286, 2.332058906555176),
    @(287, 2.49203085899353),
    @(288, 2.147327661514282),
    @(289, 2.002383232116699),
    @(290, 1.550680875778198),
    @(291, 1.520063519477844),
    @(292, 1.853135585784912),
    @(293, 1.991788864135742),
    @(294, 1.793512463569641),
    @(295, 1.649678945541382),
    @(296, 1.779754638671875),
    @(297, 1.746778845787048),
    @(298, 1.716609239578247),
    @(299, 1.769177317619324),
    @(300, 2.062208414077759),
    @(301, 2.313988447189331),
    @(302, 2.500409603118896),
    @(303, 2.339307069778442),
    @(304, 1.542866349220276)
)

for ($i = 0; $i -lt $updatedValues.Length; $i++) {
    $rowNum = $updatedValues[$i][0]
    $newVal = $updatedValues[$i][1]
    $ws.Cells.Item($rowNum, 2).Value = $newVal
}

# ---------------------------------------------------------------------
# 2) Append new rows 305-340 with date (col A) and estimate (col B)
# ---------------------------------------------------------------------
$newRows = @(
    @(43922, 0.3452092409133911),
    @(43952, 0.2264103293418884),
    @(43983, 0.7160224914550781),
    @(44013, 1.014145374298096),
    @(44044, 1.309081315994263),
    @(44075, 1.371483325958252),
    @(44105, 1.182543754577637),
    @(44136, 1.167550086975098),
    @(44166, 1.322039723396301),
    @(44197, 1.394784092903137),
    @(44228, 1.693364977836609),
    @(44256, 2.630522012710571),
    @(44287, 4.130548000335693),
    @(44317, 4.915035724639893),
    @(44348, 5.281610012054443),
    @(44378, 5.221502304077148),
    @(44409, 5.188285827636719),
    @(44440, 5.383639335632324),
    @(44470, 6.23775053024292),
    @(44501, 6.862392425537109),
    @(44531, 7.194454669952393),
    @(44562, 7.59528112411499),
    @(44593, 7.954841613769531),
    @(44621, 8.515210151672363),
    @(44652, 8.227764129638672),
    @(44682, 8.502333641052246),
    @(44713, 8.932989120483398),
    @(44743, 8.413175582885742),
    @(44774, 8.227365493774414),
    @(44805, 8.214848518371582),
    @(44835, 7.762491226196289),
    @(44866, 7.135345935821533),
    @(44896, 6.444939613342285),
    @(44927, 6.347160339355469),
    @(44958, 5.986446857452393),
    @(44986, 4.986930847167969)
)

$lastDataRow = 304
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowNum = $lastDataRow + $i + 1
    $dateSerial = $newRows[$i][0]
    $val = $newRows[$i][1]

    # Clone the formatting (style, number format, borders, etc.) of the
    # last existing date cell so the new date cell matches the rest of
    # the column, then overwrite the copied value.
    $ws.Range("A" + $lastDataRow).Copy($ws.Cells.Item($rowNum, 1))
    $ws.Cells.Item($rowNum, 1).Value = $dateSerial
    $ws.Cells.Item($rowNum, 2).Value = $val
}
